$wb = $excel.ActiveWorkbook

$wsIT = $wb.Worksheets.Item("IT Department")
$wsHR = $wb.Worksheets.Item("HR Department")

# --- Payroll Period update (both sheets share the same shared string) ---
$wsIT.Range("B1").Value = "2014-12-10-2014-12-25"
$wsHR.Range("B1").Value = "2014-12-10-2014-12-25"

# --- IT Department, row 8 (Kier Pogi Boromeo) fixes ---
# ALLOWANCE: 0 -> 1 (plain numeric cell)
$wsIT.Range("H8").Value = 1

# Use a same-style donor cell (E8) to copy number formatting onto the
# cells we rewrite as text, so the quote-prefix trick used to force text
# type doesn't leave a brand-new / different style behind.
$styleDonor = $wsIT.Range("E8")

# TOTAL ALLOWANCES: "6,000.00" -> "6,001.00" (kept as text)
$i8 = $wsIT.Range("I8")
$i8.Value = "'6,001.00"
$styleDonor.Copy()
$i8.PasteSpecial(-4122)

# NET: "5,340.75" -> "5,341.75" (kept as text)
$r8 = $wsIT.Range("R8")
$r8.Value = "'5,341.75"
$styleDonor.Copy()
$r8.PasteSpecial(-4122)

$excel.CutCopyMode = 0
